# Update cryptocurrency price/volume data to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.019.60"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").Value = "1.870.41"
$ws.Range("E3").Value = "  -2.64%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'319.52"

$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D8").Value = "'0.3934"
$ws.Range("E8").Value = "  -2.33%  "

$ws.Range("D9").Value = "'0.08204"
$ws.Range("E9").Value = "  -3.05%  "

$ws.Range("D10").Value = "'42.25"
$ws.Range("E10").Value = "  -1.23%  "

$ws.Range("D11").Value = "'1.092"
$ws.Range("E11").Value = "  -2.93%  "

$ws.Range("D12").Value = "'22.89"
$ws.Range("E12").Value = "  +1.66%  "

$ws.Range("D13").Value = "1.850.61"
$ws.Range("E13").Value = "  -3.92%  "

$ws.Range("D14").Value = "'6.271"
$ws.Range("E14").Value = "  -1.29%  "

$ws.Range("D15").Value = "'7.156"
$ws.Range("E15").Value = "  -2.87%  "

$ws.Range("E16").Value = "  +0.08%  "

$ws.Range("D17").Value = "'91.83"
$ws.Range("E17").Value = "  -4.32%  "

$ws.Range("D18").Value = "'0.00001082"
$ws.Range("E18").Value = "  -3.13%  "

$ws.Range("D19").Value = "'0.06419"
$ws.Range("E19").Value = "  -4.53%  "

$ws.Range("D20").Value = "'17.97"
$ws.Range("E20").Value = "  -1.83%  "

$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").Value = "29.999.31"
$ws.Range("E22").Value = "  -0.42%  "

$ws.Range("D23").Value = "'5.808"

$ws.Range("D24").Value = "'11.10"
$ws.Range("E24").Value = "  -1.92%  "

$ws.Range("D25").Value = "'2.165"
$ws.Range("E25").Value = "  -1.42%  "

$ws.Range("D26").Value = "2.078.11"
$ws.Range("E26").Value = "  -3.20%  "

$ws.Range("D27").Value = "'161.16"
$ws.Range("E27").Value = "  +0.24%  "

$ws.Range("D28").Value = "'20.96"
$ws.Range("E28").Value = "  -1.62%  "

$ws.Range("D29").Value = "'2.232"
$ws.Range("E29").Value = "  -9.33%  "

$ws.Range("D30").Value = "'127.07"
$ws.Range("E30").Value = "  -1.51%  "

$ws.Range("D31").Value = "'1.054"
$ws.Range("E31").Value = "  -2.56%  "

$ws.Range("D32").Value = "'0.1035"
$ws.Range("E32").Value = "  -2.29%  "

$ws.Range("D33").Value = "'5.910"
$ws.Range("E33").Value = "  -2.82%  "

$ws.Range("D34").Value = "'3.755"
$ws.Range("E34").Value = "  +2.64%  "

$ws.Range("D35").Value = "'0.02416"
$ws.Range("E35").Value = "  -4.11%  "

$ws.Range("D36").Value = "'5.275"
$ws.Range("E36").Value = "  +0.94%  "

$ws.Range("D37").Value = "'0.06331"
$ws.Range("E37").Value = "  -4.10%  "

$ws.Range("D38").Value = "'0.2144"
$ws.Range("E38").Value = "  -3.15%  "

$ws.Range("D39").Value = "'1.173"
$ws.Range("E39").Value = "  -5.35%  "

$ws.Range("D40").Value = "'8.509"
$ws.Range("E40").Value = "  -5.78%  "

$ws.Range("D41").Value = "'0.6302"
$ws.Range("E41").Value = "  -4.11%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'11.26"
$ws.Range("E42").Value = "  -3.48%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.203"
$ws.Range("E43").Value = "  -3.66%  "

$ws.Range("D44").Value = "'0.9997"
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.5897"
$ws.Range("E45").Value = "  -4.43%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.86"
$ws.Range("E46").Value = "  -2.97%  "

$ws.Range("D47").Value = "'3.629"
$ws.Range("E47").Value = "  -3.44%  "

$ws.Range("D48").Value = "'1.991"
$ws.Range("E48").Value = "  -3.60%  "

$ws.Range("D49").Value = "'122.56"
$ws.Range("E49").Value = "  -2.54%  "

$ws.Range("D50").Value = "'1.203"
$ws.Range("E50").Value = "  -3.44%  "

$ws.Range("D51").Value = "'1.141"
$ws.Range("E51").Value = "  -1.71%  "
